# Generate Report for Handoff
#
# "b.md" has now been handed off to localization for both the zh-cn and
# de-de locales, so the status/handoff columns for that row need to move
# from "Handed back: in sync with en-US" to "Ready for handoff" on the
# Overview sheet, and the per-locale sheets need their b.md row updated
# with the new handoff file name + handoff datetime.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row for b.md (row 3)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-21 14:33:55"

# ---------------------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-21 14:33:51"

# The D3 hyperlink needs its display text updated to the new handoff file
# name, keeping the same target address. This host models every hyperlink
# property write as "append a new hyperlink", so drop the stale one first.
$zhcnOldLink = $zhcn.Range("D3").Hyperlinks.Item(1)
$zhcnTarget = $zhcnOldLink.Address
$zhcnOldLink.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), $zhcnTarget, $null, $null, "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")

# ---------------------------------------------------------------------------
# de-de sheet: row for b.md (row 3)
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-21 14:33:55"

$dedeOldLink = $dede.Range("D3").Hyperlinks.Item(1)
$dedeTarget = $dedeOldLink.Address
$dedeOldLink.Delete()
$dede.Hyperlinks.Add($dede.Range("D3"), $dedeTarget, $null, $null, "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
